# Apply updated cryptocurrency price/volume data (and two row re-orderings)
# to Sheet1 of the workbook, matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values are plain text in the source sheet (e.g. "26.934.20" or
# "0.4820"); assigning a leading apostrophe forces Excel to keep them as text instead
# of auto-converting to numbers (which would drop trailing zeros / change formatting).
# Row 2
$ws.Range("D2").Value = "'" + '26.934.20'
$ws.Range("E2").Value = '  -0.62%  '
# Row 3
$ws.Range("D3").Value = "'" + '1.821.42'
$ws.Range("E3").Value = '  -0.15%  '
# Row 4
$ws.Range("D4").Value = "'" + '1.003'
$ws.Range("E4").Value = '  -0.64%  '
# Row 5
$ws.Range("D5").Value = "'" + '310.71'
$ws.Range("E5").Value = '  -0.27%  '
# Row 6
$ws.Range("D6").Value = "'" + '1.003'
$ws.Range("E6").Value = '  -0.53%  '
# Row 7
$ws.Range("D7").Value = "'" + '0.4622'
$ws.Range("E7").Value = '  -0.95%  '
# Row 8
$ws.Range("D8").Value = "'" + '0.3702'
$ws.Range("E8").Value = '  +1.64%  '
# Row 9
$ws.Range("D9").Value = "'" + '0.07329'
$ws.Range("E9").Value = '  +0.39%  '
# Row 10
$ws.Range("D10").Value = "'" + '0.8730'
$ws.Range("E10").Value = '  +0.45%  '
# Row 11
$ws.Range("D11").Value = "'" + '0.07890'
$ws.Range("E11").Value = '  +3.52%  '
# Row 12
$ws.Range("D12").Value = "'" + '19.68'
$ws.Range("E12").Value = '  -2.49%  '
# Row 13
$ws.Range("D13").Value = "'" + '1.882.01'
$ws.Range("E13").Value = '  +1.25%  '
# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'" + '5.325'
$ws.Range("E14").Value = '  -0.33%  '
# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = "'" + '6.542'
$ws.Range("E15").Value = '  +1.03%  '
# Row 16
$ws.Range("D16").Value = "'" + '91.24'
$ws.Range("E16").Value = '  -2.00%  '
# Row 17
$ws.Range("D17").Value = "'" + '1.005'
$ws.Range("E17").Value = '  -0.41%  '
# Row 18
$ws.Range("D18").Value = "'" + '0.000008823'
$ws.Range("E18").Value = '  +2.08%  '
# Row 19
$ws.Range("D19").Value = "'" + '1.004'
$ws.Range("E19").Value = '  -0.43%  '
# Row 20
$ws.Range("E20").Value = '  +1.89%  '
# Row 21
$ws.Range("D21").Value = "'" + '26.966.06'
$ws.Range("E21").Value = '  -1.01%  '
# Row 22
$ws.Range("D22").Value = "'" + '5.093'
$ws.Range("E22").Value = '  -1.89%  '
# Row 23
$ws.Range("E23").Value = '  -0.53%  '
# Row 24
$ws.Range("D24").Value = "'" + '2.077.41'
$ws.Range("E24").Value = '  +0.26%  '
# Row 25
$ws.Range("E25").Value = '  +1.15%  '
# Row 26
$ws.Range("D26").Value = "'" + '1.848'
$ws.Range("E26").Value = '  -1.16%  '
# Row 27
$ws.Range("D27").Value = "'" + '18.34'
$ws.Range("E27").Value = '  +0.38%  '
# Row 28
$ws.Range("D28").Value = "'" + '2.030'
$ws.Range("E28").Value = '  -3.51%  '
# Row 29
$ws.Range("E29").Value = '  +0.53%  '
# Row 30
$ws.Range("D30").Value = "'" + '115.54'
$ws.Range("E30").Value = '  -0.53%  '
# Row 31
$ws.Range("D31").Value = "'" + '0.08871'
$ws.Range("E31").Value = '  -0.63%  '
# Row 32
$ws.Range("D32").Value = "'" + '2.957'
$ws.Range("E32").Value = '  -0.04%  '
# Row 33
$ws.Range("D33").Value = "'" + '0.7281'
$ws.Range("E33").Value = '  -0.68%  '
# Row 34
$ws.Range("D34").Value = "'" + '4.431'
$ws.Range("E34").Value = '  -0.61%  '
# Row 35
$ws.Range("E35").Value = '  -1.19%  '
# Row 36
$ws.Range("D36").Value = "'" + '2.460'
$ws.Range("E36").Value = '  -3.08%  '
# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'" + '0.01944'
$ws.Range("E37").Value = '  +1.27%  '
# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = "'" + '1.067'
$ws.Range("E38").Value = '  -0.44%  '
# Row 39
$ws.Range("D39").Value = "'" + '0.05209'
$ws.Range("E39").Value = '  -1.10%  '
# Row 40
$ws.Range("D40").Value = "'" + '2.947'
$ws.Range("E40").Value = '  +0.32%  '
# Row 41
$ws.Range("D41").Value = "'" + '7.075'
$ws.Range("E41").Value = '  -0.67%  '
# Row 42
$ws.Range("D42").Value = "'" + '0.5144'
$ws.Range("E42").Value = '  -1.55%  '
# Row 43
$ws.Range("E43").Value = '  -0.70%  '
# Row 44
$ws.Range("D44").Value = "'" + '8.143'
$ws.Range("E44").Value = '  -1.57%  '
# Row 45
$ws.Range("D45").Value = "'" + '0.4820'
$ws.Range("E45").Value = '  -1.04%  '
# Row 46
$ws.Range("D46").Value = "'" + '1.003'
$ws.Range("E46").Value = '  -0.60%  '
# Row 47
$ws.Range("D47").Value = "'" + '10.16'
$ws.Range("E47").Value = '  +0.31%  '
# Row 48
$ws.Range("D48").Value = "'" + '102.51'
$ws.Range("E48").Value = '  -1.18%  '
# Row 49
$ws.Range("D49").Value = "'" + '1.628'
$ws.Range("E49").Value = '  -0.60%  '
# Row 50
$ws.Range("D50").Value = "'" + '0.06186'
$ws.Range("E50").Value = '  -0.96%  '
# Row 51
$ws.Range("D51").Value = "'" + '64.89'
$ws.Range("E51").Value = '  +0.54%  '
